$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("M2").Value = 11.05182166666667
$ws.Range("N2").Value = 33.155465
$ws.Range("O2").Value = 0.09655996768538078
$ws.Range("P2").Value = 0.09655996768538078
$ws.Range("Q2").Value = 1877.877539538502
$ws.Range("R2").Value = 16900.89785584652
$ws.Range("S2").Value = 0.04288406734930842
$ws.Range("T2").Value = 0.04288406734930843
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("O3").Value = 0.620602129029037
$ws.Range("P3").Value = 0.620602129029037
$ws.Range("Q3").Value = 12069.33708688314
$ws.Range("R3").Value = 108624.0337819483
$ws.Range("S3").Value = 0.2756208824046115
$ws.Range("T3").Value = 0.2756208824046115
$ws.Range("G4").Value = 169.915657
$ws.Range("H4").Value = 509.746971
$ws.Range("I4").Value = 0.4441184931734509
$ws.Range("J4").Value = 0.4441184931734509
$ws.Range("M4").Value = 32.37236033333333
$ws.Range("N4").Value = 97.117081
$ws.Range("O4").Value = 0.2828379032855822
$ws.Range("P4").Value = 0.2828379032855822
$ws.Range("Q4").Value = 5500.570874679072
$ws.Range("R4").Value = 49505.13787211166
$ws.Range("S4").Value = 0.125613543419531
$ws.Range("T4").Value = 0.125613543419531
$ws.Range("I5").Value = 0.1787346690539575
$ws.Range("J5").Value = 0.1787346690539575
$ws.Range("M5").Value = 11.05182166666667
$ws.Range("N5").Value = 33.155465
$ws.Range("O5").Value = 0.09655996768538078
$ws.Range("P5").Value = 0.09655996768538078
$ws.Range("Q5").Value = 755.7483548026651
$ws.Range("R5").Value = 6801.735193223984
$ws.Range("S5").Value = 0.01725861386810737
$ws.Range("T5").Value = 0.01725861386810737
$ws.Range("I6").Value = 0.1787346690539575
$ws.Range("J6").Value = 0.1787346690539575
$ws.Range("O6").Value = 0.620602129029037
$ws.Range("P6").Value = 0.620602129029037
$ws.Range("S6").Value = 0.1109231161461864
$ws.Range("T6").Value = 0.1109231161461864
$ws.Range("I7").Value = 0.1787346690539575
$ws.Range("J7").Value = 0.1787346690539575
$ws.Range("M7").Value = 32.37236033333333
$ws.Range("N7").Value = 97.117081
$ws.Range("O7").Value = 0.2828379032855822
$ws.Range("P7").Value = 0.2828379032855822
$ws.Range("Q7").Value = 2213.694610797561
$ws.Range("R7").Value = 19923.25149717805
$ws.Range("S7").Value = 0.05055293903966378
$ws.Range("T7").Value = 0.05055293903966378
$ws.Range("G8").Value = 53.27463399999999
$ws.Range("H8").Value = 159.823902
$ws.Range("I8").Value = 0.1392470275793777
$ws.Range("J8").Value = 0.1392470275793778
$ws.Range("M8").Value = 11.05182166666667
$ws.Range("N8").Value = 33.155465
$ws.Range("O8").Value = 0.09655996768538078
$ws.Range("P8").Value = 0.09655996768538078
$ws.Range("Q8").Value = 588.7817543249366
$ws.Range("R8").Value = 5299.035788924429
$ws.Range("S8").Value = 0.01344568848335004
$ws.Range("T8").Value = 0.01344568848335004
$ws.Range("G9").Value = 53.27463399999999
$ws.Range("H9").Value = 159.823902
$ws.Range("I9").Value = 0.1392470275793777
$ws.Range("J9").Value = 0.1392470275793778
$ws.Range("O9").Value = 0.620602129029037
$ws.Range("P9").Value = 0.620602129029037
$ws.Range("Q9").Value = 3784.168729820616
$ws.Range("R9").Value = 34057.51856838555
$ws.Range("S9").Value = 0.08641700177672686
$ws.Range("T9").Value = 0.08641700177672687
$ws.Range("G10").Value = 53.27463399999999
$ws.Range("H10").Value = 159.823902
$ws.Range("I10").Value = 0.1392470275793777
$ws.Range("J10").Value = 0.1392470275793778
$ws.Range("M10").Value = 32.37236033333333
$ws.Range("N10").Value = 97.117081
$ws.Range("O10").Value = 0.2828379032855822
$ws.Range("P10").Value = 0.2828379032855822
$ws.Range("Q10").Value = 1724.625648474451
$ws.Range("R10").Value = 15521.63083627006
$ws.Range("S10").Value = 0.03938433731930084
$ws.Range("T10").Value = 0.03938433731930085
$ws.Range("G11").Value = 91.01828266666666
$ws.Range("H11").Value = 273.054848
$ws.Range("I11").Value = 0.2378998101932138
$ws.Range("J11").Value = 0.2378998101932138
$ws.Range("M11").Value = 11.05182166666667
$ws.Range("N11").Value = 33.155465
$ws.Range("O11").Value = 0.09655996768538078
$ws.Range("P11").Value = 0.09655996768538078
$ws.Range("Q11").Value = 1005.917828438258
$ws.Range("R11").Value = 9053.260455944319
$ws.Range("S11").Value = 0.02297159798461494
$ws.Range("T11").Value = 0.02297159798461495
$ws.Range("G12").Value = 91.01828266666666
$ws.Range("H12").Value = 273.054848
$ws.Range("I12").Value = 0.2378998101932138
$ws.Range("J12").Value = 0.2378998101932138
$ws.Range("O12").Value = 0.620602129029037
$ws.Range("P12").Value = 0.620602129029037
$ws.Range("Q12").Value = 6465.150734009244
$ws.Range("R12").Value = 58186.35660608319
$ws.Range("S12").Value = 0.1476411287015123
$ws.Range("T12").Value = 0.1476411287015123
$ws.Range("G13").Value = 91.01828266666666
$ws.Range("H13").Value = 273.054848
$ws.Range("I13").Value = 0.2378998101932138
$ws.Range("J13").Value = 0.2378998101932138
$ws.Range("M13").Value = 32.37236033333333
$ws.Range("N13").Value = 97.117081
$ws.Range("O13").Value = 0.2828379032855822
$ws.Range("P13").Value = 0.2828379032855822
$ws.Range("Q13").Value = 2946.476643406521
$ws.Range("R13").Value = 26518.28979065869
$ws.Range("S13").Value = 0.06728708350708658
$ws.Range("T13").Value = 0.06728708350708658
